# Update the "volume" column (D) data values on the active sheet
# and move the active selection, per the commit "finish bildings, add spy_fact_temp".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 31914
$ws.Range("D3").Value = 36167
$ws.Range("D4").Value = 45221
$ws.Range("D5").Value = 32960
$ws.Range("D6").Value = 118910
$ws.Range("D7").Value = 86979
$ws.Range("D8").Value = 17105
$ws.Range("D9").Value = 517516
$ws.Range("D10").Value = 469399

$ws.Range("D11").Select()
